$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.173.30"
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").Value = "1.928.86"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'330.82"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "'0.4732"
$ws.Range("E7").Value = "  -4.80%  "
$ws.Range("D8").Value = "'0.4059"
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("D9").Value = "'53.01"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.08465"
$ws.Range("E10").Value = "  -8.07%  "
$ws.Range("D11").Value = "'1.049"
$ws.Range("E11").Value = "  -4.33%  "
$ws.Range("D12").Value = "'22.26"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("D13").Value = "1.942.23"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "'7.531"
$ws.Range("E14").Value = "  -4.40%  "
$ws.Range("D15").Value = "'6.107"
$ws.Range("E15").Value = "  -5.45%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "'90.25"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").Value = "'0.06592"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "'18.22"
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'5.786"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").Value = "28.212.39"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "'11.44"
$ws.Range("E24").Value = "  -4.68%  "
$ws.Range("D25").Value = "'2.286"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "2.201.31"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'154.62"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").Value = "'20.13"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "'2.158"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").Value = "'5.750"
$ws.Range("E30").Value = "  -9.87%  "
$ws.Range("D31").Value = "'123.93"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'0.9813"
$ws.Range("E32").Value = "  -6.08%  "
$ws.Range("D33").Value = "'0.09617"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "'1.447"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").Value = "'5.579"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").Value = "'3.643"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").Value = "'9.132"
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").Value = "'0.02324"
$ws.Range("E38").Value = "  -4.44%  "
$ws.Range("D39").Value = "'0.06180"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").Value = "'1.240"
$ws.Range("E40").Value = "  -5.93%  "
$ws.Range("D41").Value = "'0.6185"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "'1.005"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "'0.1904"
$ws.Range("E44").Value = "  -4.06%  "
$ws.Range("D45").Value = "'1.316"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "'0.5894"
$ws.Range("E46").Value = "  -5.05%  "
$ws.Range("D47").Value = "'12.79"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "'2.040"
$ws.Range("E48").Value = "  -7.12%  "
$ws.Range("D49").Value = "'3.471"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'0.06804"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "'109.92"
$ws.Range("E51").Value = "  -2.65%  "
